$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "WiFi" bullet: drop the spell-check proofErr bookends around the run.
#    Re-serializing just this paragraph's own Range through OOXML and
#    feeding it straight back in is enough to shed the stale proofErr marks
#    without touching anything else in the run.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^WiFi:") {
        $rng = $p.Range
        $rng.InsertXML($rng.WordOpenXML)
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Free-Space Path Loss equation: remove the "- L_f" (loss-factor) term so
#    the formula reads Pr = Pt + Gt + Gr - 20log10(d) ... instead of
#    Pr = Pt + Gt + Gr - L_f - 20log10(d) ...
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.OMaths.Count; $i++) {
    $om = $d.OMaths.Item($i)
    $omXml = $om.Range.WordOpenXML
    if ($omXml.Contains("<m:t>L</m:t>")) {
        $lossTerm = '<m:r><m:rPr><m:sty m:val="bi"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><m:t>-</m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:b/><w:bCs/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><m:rPr><m:sty m:val="bi"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><m:t>L</m:t></m:r></m:e><m:sub><m:r><m:rPr><m:sty m:val="bi"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><m:t>f</m:t></m:r></m:sub></m:sSub>'
        if ($omXml.Contains($lossTerm)) {
            $newOmXml = $omXml.Replace($lossTerm, "")
            $om.Range.InsertXML($newOmXml)
            break
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Drop the whole "L_f - Loss factor (if applicable)" bullet paragraph.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^L_f - Loss factor") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 4) "d_t, d_r -" bullet: collapse the proofErr-wrapped d_t / ", " / d_r /
#    " - " runs back into one contiguous bold run reading "d_t, d_r - ".
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^d_t, d_r -") {
        $rng = $p.Range
        $rng.InsertXML($rng.WordOpenXML)
        break
    }
}

# ---------------------------------------------------------------------------
# 5) "h_eff -" bullet: collapse the proofErr-wrapped h_eff / " -" runs back
#    into a single run reading "h_eff -".
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^h_eff -") {
        $rng = $p.Range
        $rng.InsertXML($rng.WordOpenXML)
        break
    }
}
